$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 250004750
$ws.Range("I76").Value = 500002000
$ws.Range("J76").Value = 7499.5
$ws.Range("K76").Value = 500002000
$ws.Range("L76").Value = 7499.5
$ws.Range("M76").Value = -500001685
$ws.Range("N76").Value = -8129.5
$ws.Range("H79").Value = 250004750
$ws.Range("I79").Value = 500002000
$ws.Range("J79").Value = 7499.5
$ws.Range("K79").Value = 500002000
$ws.Range("L79").Value = 7499.5
$ws.Range("M79").Value = -500000908
$ws.Range("N79").Value = -9683.5
$ws.Range("H88").Value = 3165.125
$ws.Range("I88").Value = 1286.3334
$ws.Range("J88").Value = 4292.4
$ws.Range("K88").Value = 1286.3334
$ws.Range("L88").Value = 4292.4
$ws.Range("M88").Value = -880.3334
$ws.Range("N88").Value = -5104.4
$ws.Range("H91").Value = 3165.125
$ws.Range("I91").Value = 1286.3334
$ws.Range("J91").Value = 4292.4
$ws.Range("K91").Value = 1286.3334
$ws.Range("L91").Value = 4292.4
$ws.Range("M91").Value = 117.6666
$ws.Range("N91").Value = -7100.4
$ws.Range("H100").Value = 3223.8096
$ws.Range("I100").Value = 3117.7646
$ws.Range("K100").Value = 3117.7646
$ws.Range("M100").Value = -2576.7646
$ws.Range("H135").Value = 2491.5625
$ws.Range("I135").Value = 1104.925
$ws.Range("J135").Value = 9424.75
$ws.Range("K135").Value = 9944.324999999999
$ws.Range("L135").Value = 84822.75
$ws.Range("M135").Value = -7409.324999999999
$ws.Range("N135").Value = -89892.75
$ws.Range("H138").Value = 4201.731
$ws.Range("I138").Value = 1732.238
$ws.Range("K138").Value = 5196.714
$ws.Range("M138").Value = -56.71399999999994
$ws.Range("H141").Value = 693.9
$ws.Range("I141").Value = 493.22223
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 1479.66669
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = 3700.33331
$ws.Range("N141").Value = -17860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1823661.4
$ws.Range("I32").Value = 1350.3636
$ws.Range("K32").Value = 1350.3636
$ws.Range("M32").Value = -1063.3636
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H124").Value = 31247.75
$ws.Range("J124").Value = 31247.75
$ws.Range("L124").Value = 31247.75
$ws.Range("N124").Value = -41067.75
$ws.Range("H132").Value = 981414.8
$ws.Range("I132").Value = 1264748.9
$ws.Range("J132").Value = 131412.5
$ws.Range("K132").Value = 3794246.7
$ws.Range("L132").Value = 394237.5
$ws.Range("M132").Value = -3791716.7
$ws.Range("N132").Value = -399297.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8334363
$ws.Range("I107").Value = 12500393
$ws.Range("K107").Value = 12500393
$ws.Range("M107").Value = -12498473
$ws.Range("H124").Value = 74999.5
$ws.Range("J124").Value = 74999.5
$ws.Range("L124").Value = 74999.5
$ws.Range("N124").Value = -84819.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3750.5264
$ws.Range("I31").Value = 1830.125
$ws.Range("J31").Value = 5147.1816
$ws.Range("K31").Value = 1830.125
$ws.Range("L31").Value = 5147.1816
$ws.Range("M31").Value = -1535.125
$ws.Range("N31").Value = -5737.1816
$ws.Range("H34").Value = 3750.5264
$ws.Range("I34").Value = 1830.125
$ws.Range("J34").Value = 5147.1816
$ws.Range("K34").Value = 1830.125
$ws.Range("L34").Value = 5147.1816
$ws.Range("M34").Value = -1628.125
$ws.Range("N34").Value = -5551.1816
$ws.Range("H107").Value = 833.125
$ws.Range("I107").Value = 681.3333
$ws.Range("J107").Value = 924.2
$ws.Range("K107").Value = 681.3333
$ws.Range("L107").Value = 924.2
$ws.Range("M107").Value = 1238.6667
$ws.Range("N107").Value = -4764.2
$ws.Range("H122").Value = 1808.2273
$ws.Range("I122").Value = 1786.3334
$ws.Range("J122").Value = 1855.1428
$ws.Range("K122").Value = 5359.0002
$ws.Range("L122").Value = 5565.428400000001
$ws.Range("M122").Value = -2909.0002
$ws.Range("N122").Value = -10465.4284
$ws.Range("H132").Value = 4229.5
$ws.Range("I132").Value = 3549.7083
$ws.Range("K132").Value = 10649.1249
$ws.Range("M132").Value = -8119.124899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7333601
$ws.Range("I4").Value = 7333601
$ws.Range("K4").Value = 22000803
$ws.Range("M4").Value = -22000691
$ws.Range("H26").Value = 340.7
$ws.Range("I26").Value = 382.81818
$ws.Range("K26").Value = 1148.45454
$ws.Range("M26").Value = -860.45454
$ws.Range("H113").Value = 313282.7
$ws.Range("J113").Value = 357983.5
$ws.Range("L113").Value = 1073950.5
$ws.Range("N113").Value = -1078290.5
$ws.Range("H121").Value = 29529.908
$ws.Range("I121").Value = 873.75
$ws.Range("J121").Value = 45904.855
$ws.Range("K121").Value = 2621.25
$ws.Range("L121").Value = 137714.565
$ws.Range("M121").Value = -1311.25
$ws.Range("N121").Value = -140334.565
$ws.Range("H131").Value = 50984816
$ws.Range("J131").Value = 66669868
$ws.Range("L131").Value = 200009604
$ws.Range("N131").Value = -200019684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28040
$ws.Range("I70").Value = 10050
$ws.Range("K70").Value = 10050
$ws.Range("M70").Value = -9780
$ws.Range("H73").Value = 28040
$ws.Range("I73").Value = 10050
$ws.Range("K73").Value = 10050
$ws.Range("M73").Value = -9114

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 5016299.5
$ws.Range("I11").Value = 5016299.5
$ws.Range("K11").Value = 5016299.5
$ws.Range("M11").Value = -5016159.5
$ws.Range("H16").Value = 2807.075
$ws.Range("I16").Value = 846.67645
$ws.Range("J16").Value = 13916
$ws.Range("K16").Value = 846.67645
$ws.Range("L16").Value = 13916
$ws.Range("M16").Value = -676.67645
$ws.Range("N16").Value = -14256
$ws.Range("H22").Value = 796.5
$ws.Range("I22").Value = 733.1429000000001
$ws.Range("K22").Value = 733.1429000000001
$ws.Range("M22").Value = -438.1429000000001
$ws.Range("H25").Value = 8358943.5
$ws.Range("I25").Value = 8358943.5
$ws.Range("K25").Value = 8358943.5
$ws.Range("M25").Value = -8358713.5
$ws.Range("H27").Value = 796.5
$ws.Range("I27").Value = 733.1429000000001
$ws.Range("K27").Value = 733.1429000000001
$ws.Range("M27").Value = -626.1429000000001
$ws.Range("H61").Value = 5158.2163
$ws.Range("I61").Value = 4647.6
$ws.Range("J61").Value = 5758.9414
$ws.Range("K61").Value = 4647.6
$ws.Range("L61").Value = 5758.9414
$ws.Range("M61").Value = -4445.6
$ws.Range("N61").Value = -6162.9414
$ws.Range("H113").Value = 5158.2163
$ws.Range("I113").Value = 4647.6
$ws.Range("J113").Value = 5758.9414
$ws.Range("K113").Value = 4647.6
$ws.Range("L113").Value = 5758.9414
$ws.Range("M113").Value = -2477.6
$ws.Range("N113").Value = -10098.9414
$ws.Range("H132").Value = 6065.5
$ws.Range("I132").Value = 5329.8184
$ws.Range("J132").Value = 7414.25
$ws.Range("K132").Value = 15989.4552
$ws.Range("L132").Value = 22242.75
$ws.Range("M132").Value = -13459.4552
$ws.Range("N132").Value = -27302.75
$ws.Range("H136").Value = 22729810
$ws.Range("I136").Value = 31251944
$ws.Range("K136").Value = 93755832
$ws.Range("M136").Value = -93753282

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1310.3235
$ws.Range("J81").Value = 1700
$ws.Range("L81").Value = 3400
$ws.Range("N81").Value = -5522
$ws.Range("H84").Value = 1310.3235
$ws.Range("J84").Value = 1700
$ws.Range("L84").Value = 17000
$ws.Range("N84").Value = -27608
$ws.Range("H107").Value = 885.0526
$ws.Range("J107").Value = 1109.4445
$ws.Range("L107").Value = 3328.3335
$ws.Range("N107").Value = -7168.333500000001
$ws.Range("H132").Value = 13081.272
$ws.Range("I132").Value = 9379.6
$ws.Range("J132").Value = 16166
$ws.Range("K132").Value = 28138.8
$ws.Range("L132").Value = 48498
$ws.Range("M132").Value = -25608.8
$ws.Range("N132").Value = -53558
$ws.Range("H136").Value = 11909380
$ws.Range("I136").Value = 16667755
$ws.Range("J136").Value = 13445.25
$ws.Range("K136").Value = 50003265
$ws.Range("M136").Value = -50000715
$ws.Range("N136").Value = -45435.75
